# Update "Pais" COVID-19 country statistics sheet with the latest refresh
# (countries & provincias Spain), matching the 4 Abr 2020 20:22 data pull.
# The source table is re-sorted by "Casos totales" (column B) descending,
# so several rows shift position in addition to the numbers themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @(4, "Estados Unidos", 300720,23559,14518,278040,7869,758,8162),
    @(5, "España", 124736,5537,34219,78773,6416,546,11744),
    @(6, "Italia", 124632,4805,20996,88274,3994,681,15362),
    @(7, "Alemania", 92150,991,26400,64420,3936,55,1330),
    @(8, "Francia", 89953,7788,15438,66955,6838,1053,7560),
    @(9, "China", 81639,19,76755,1558,331,4,3326),
    @(10, "Iran", 55743,2560,19736,32555,4103,158,3452),
    @(11, "Reino Unido", 41903,3735,135,37455,163,708,4313),
    @(12, "Turquia", 23934,3013,786,22647,1311,76,501),
    @(13, "Suiza", 20278,672,4846,14791,391,50,641),
    @(14, "Belgica", 18431,1661,3247,13901,1245,140,1283),
    @(15, "Paises Bajos", 16627,904,250,14726,1360,164,1651),
    @(16, "Canada", 13860,1485,2559,11073,120,20,228),
    @(17, "Austria", 11781,257,2507,9088,245,18,186),
    @(18, "Portugal", 10524,638,75,10183,251,20,266),
    @(19, "Corea del Sur", 10156,94,6325,3654,55,3,177),
    @(20, "Brasil", 9391,197,127,8888,296,13,376),
    @(21, "Israel", 7851,423,427,7380,115,4,44),
    @(22, "Suecia", 6443,312,205,5865,379,15,373),
    @(23, "Noruega", 5550,180,32,5456,98,3,62),
    @(24, "Australia", 5550,96,585,4935,85,2,30),
    @(25, "Rusia", 4731,582,333,4355,8,9,43),
    @(26, "Irlanda", 4604,331,25,4442,148,17,137),
    @(27, "Chequia", 4362,172,78,4225,87,6,59),
    @(28, "Chile", 4161,424,528,3606,38,5,27),
    @(29, "Dinamarca", 4077,320,1283,2633,142,22,161),
    @(30, "Polonia", 3627,244,116,3432,50,8,79),
    @(31, "Rumania", 3613,430,329,3138,119,13,146),
    @(32, "Malasia", 3483,150,915,2511,99,4,57),
    @(33, "Ecuador", 3465,97,65,3228,100,27,172),
    @(34, "Filipinas", 3094,76,57,2893,1,8,144),
    @(35, "India", 3082,23,229,2767,0,0,86),
    @(36, "Japon", 2935,0,514,2352,60,0,69),
    @(37, "Pakistan", 2818,132,131,2646,13,1,41),
    @(38, "Luxemburgo", 2729,117,500,2198,35,0,31),
    @(39, "Arabia Saudita", 2179,140,420,1730,41,4,29),
    @(40, "Indonesia", 2092,106,150,1751,0,10,191),
    @(41, "Tailandia", 2067,89,674,1373,23,1,20),
    @(42, "Finlandia", 1882,267,300,1557,72,5,25),
    @(43, "Peru", 1746,151,799,886,88,0,61),
    @(44, "Mexico", 1688,178,633,995,1,10,60),
    @(45, "Panama", 1673,0,13,1619,50,0,41),
    @(46, "Grecia", 1673,60,78,1527,92,5,68),
    @(47, "Serbia", 1624,148,54,1526,89,5,44),
    @(48, "Sudafrica", 1505,0,95,1401,7,0,9),
    @(49, "Emiratos Arabes Unidos", 1505,241,125,1370,2,1,10),
    @(50, "Republica Dominicana", 1488,0,16,1404,147,0,68),
    @(51, "Islandia", 1417,53,396,1017,12,0,4),
    @(52, "Argentina", 1353,0,266,1045,0,0,42),
    @(53, "Catar", 1325,250,109,1213,37,0,3),
    @(54, "Colombia", 1267,0,55,1187,50,0,25),
    @(55, "Argelia", 1251,80,90,1031,0,25,130),
    @(56, "Singapur", 1189,75,297,886,24,1,6),
    @(57, "Croacia", 1126,47,119,995,39,4,12),
    @(58, "Ucrania", 1115,43,23,1064,16,1,28),
    @(59, "Egipto", 1070,85,241,758,0,5,71),
    @(60, "Estonia", 1039,78,59,967,20,1,13),
    @(61, "Eslovenia", 977,43,79,876,31,2,22),
    @(62, "Nueva Zelanda", 950,82,127,822,1,0,1),
    @(63, "Marruecos", 883,92,65,760,1,10,58),
    @(64, "Irak", 878,58,259,563,0,2,56),
    @(65, "Hong Kong", 862,17,173,685,8,0,4),
    @(66, "Lituania", 771,75,7,754,11,1,10),
    @(67, "Armenia", 770,34,43,720,30,0,7),
    @(68, "Moldavia", 752,161,29,711,65,4,12),
    @(69, "Crucero", 712,0,619,82,10,0,11),
    @(70, "Barein", 688,16,423,261,3,0,4),
    @(71, "Hungria", 678,55,58,588,17,6,32),
    @(72, "Bosnia y Herzegovina", 617,38,28,570,4,2,19),
    @(73, "Tunez", 553,58,5,529,26,1,19),
    @(74, "Kazajistan", 525,61,36,484,6,2,5),
    @(75, "Azerbaiyan", 521,78,32,484,17,0,5),
    @(76, "Libano", 520,12,54,449,26,0,17),
    @(77, "Letonia", 509,16,1,507,3,0,1),
    @(78, "Camerun", 509,0,17,484,0,0,8),
    @(79, "Bulgaria", 503,18,34,452,26,3,17),
    @(80, "Republica de Macedonia", 483,53,20,446,11,5,17),
    @(81, "Kuwait", 479,62,93,385,17,1,1),
    @(82, "Eslovaquia", 471,21,10,460,3,0,1),
    @(83, "Principado de Andorra", 466,27,21,428,12,1,17),
    @(84, "Bielorrusia", 440,89,53,382,11,1,5),
    @(85, "Republica de Chipre", 426,30,33,382,11,0,11),
    @(86, "Costa Rica", 416,0,11,403,13,0,2),
    @(87, "Uruguay", 386,0,86,296,13,0,4),
    @(88, "Taiwan", 355,7,50,300,0,0,5),
    @(89, "Reunion", 334,13,40,294,4,0,0),
    @(90, "Albania", 333,29,99,216,7,1,18),
    @(91, "Jordania", 323,13,74,244,5,0,5),
    @(92, "Burkina Faso", 318,16,66,236,0,0,16),
    @(93, "Afganistan", 299,18,10,282,0,1,7),
    @(94, "Cuba", 288,19,15,267,11,0,6),
    @(95, "Oman", 277,25,61,214,3,1,2),
    @(96, "Uzbekistan", 266,39,25,239,8,0,2),
    @(97, "Honduras", 264,42,3,246,10,0,15),
    @(98, "San Marino", 251,0,26,193,13,0,32),
    @(99, "Vietnam", 240,1,90,150,3,0,0),
    @(100, "Senegal", 219,12,72,145,1,1,2),
    @(101, "Costa de Marfil", 218,0,19,198,0,0,1),
    @(102, "Estado de Palestina", 216,22,21,194,0,0,1),
    @(103, "Malta", 213,11,2,211,3,0,0),
    @(104, "Nigeria", 210,0,25,181,2,0,4),
    @(105, "Ghana", 205,0,31,169,2,0,5),
    @(106, "Montenegro", 201,27,1,198,4,0,2),
    @(107, "Mauricio", 196,10,0,189,1,0,7),
    @(108, "Islas Feroe", 181,2,93,88,2,0,0),
    @(109, "Sri Lanka", 166,7,27,134,5,1,5),
    @(110, "Georgia", 159,4,31,127,6,1,1),
    @(111, "Venezuela", 155,2,52,96,6,0,7),
    @(112, "Consejo Danes para los Refugiados", 148,0,3,129,0,0,16),
    @(113, "Kirguistan", 144,14,9,134,5,0,1),
    @(114, "Martinica", 143,0,27,113,18,0,3),
    @(115, "Bolivia", 139,7,1,128,3,1,10),
    @(116, "Brunei", 135,1,66,68,3,0,1),
    @(117, "Mayotte", 134,6,14,118,3,0,2),
    @(118, "Guadalupe", 130,0,24,99,14,0,7),
    @(119, "Isla de Man", 126,12,0,125,0,0,1),
    @(120, "Kenia", 126,4,4,118,2,0,4),
    @(121, "Niger", 120,0,0,115,0,0,5),
    @(122, "Camboya", 114,0,50,64,1,0,0),
    @(123, "Trinidad yTobago", 101,3,1,94,0,0,6),
    @(124, "Gibraltar", 98,3,52,46,0,0,0),
    @(125, "Paraguay", 96,4,12,81,2,0,3),
    @(126, "Ruanda", 89,0,0,89,0,0,0),
    @(127, "Liechtenstein", 77,2,0,77,0,0,0),
    @(128, "Guinea", 73,0,2,71,0,0,0),
    @(129, "Madagascar", 70,0,0,70,6,0,0),
    @(130, "Banglades", 70,9,30,32,1,2,8),
    @(131, "Aruba", 64,2,1,63,0,0,0),
    @(132, "Monaco", 64,0,3,60,2,0,1),
    @(133, "Guatemala", 57,7,15,40,1,1,2),
    @(134, "Guayana Francesa", 57,0,22,35,1,0,0),
    @(135, "El Salvador", 56,10,0,53,4,0,3),
    @(136, "Jamaica", 53,0,7,43,0,0,3),
    @(137, "Barbados", 52,1,0,52,0,0,0),
    @(138, "Republica de Yibuti", 50,1,8,42,0,0,0),
    @(139, "Uganda", 48,0,0,48,0,0,0),
    @(140, "Macao", 43,1,10,33,0,0,0),
    @(141, "Mali", 41,2,1,37,0,0,3),
    @(142, "Togo", 41,1,17,21,0,0,3),
    @(143, "Polinesia Francesa", 39,0,0,39,1,0,0),
    @(144, "Zambia", 39,0,2,36,0,0,1),
    @(145, "Puerto Rico", 39,0,1,36,0,0,2),
    @(146, "Etiopia", 38,3,4,34,1,0,0),
    @(147, "Bermudas", 35,0,14,21,0,0,0),
    @(148, "Guam", 32,0,0,31,0,0,1),
    @(149, "Eritrea", 29,7,0,29,0,0,0),
    @(150, "Islas Caimanes", 29,0,1,27,0,0,1),
    @(151, "Bahamas", 24,0,0,21,1,0,3),
    @(152, "San Martin (Parte Francesa)", 24,2,5,17,6,1,2),
    @(153, "Guyana", 23,0,0,19,0,0,4),
    @(154, "San Martin (Parte Holandesa)", 23,0,6,13,0,0,4),
    @(155, "Congo", 22,0,2,18,0,0,2),
    @(156, "Birmania", 21,1,0,20,0,0,1),
    @(157, "Gabon", 21,0,1,19,0,0,1),
    @(158, "Haiti", 20,2,1,19,0,0,0),
    @(159, "Tanzania", 20,0,3,16,0,0,1),
    @(160, "Maldivas", 19,0,13,6,0,0,0),
    @(161, "Guinea-Bisau", 18,3,0,18,0,0,0),
    @(162, "Nueva Caledonia", 18,0,1,17,0,0,0),
    @(163, "Islas Virgenes de los Estados Unidos", 17,0,0,17,0,0,0),
    @(164, "Libia", 17,0,0,16,0,0,1),
    @(165, "Guinea Ecuatorial", 16,0,1,15,0,0,0),
    @(166, "Benin", 16,0,2,14,0,0,0),
    @(167, "Siria", 16,0,2,12,0,0,2),
    @(168, "Antigua y Barbuda", 15,0,0,15,1,0,0),
    @(169, "Dominica", 14,0,0,14,0,0,0),
    @(170, "Mongolia", 14,0,2,12,0,0,0),
    @(171, "Namibia", 14,0,3,11,0,0,0),
    @(172, "Santa Lucia", 13,0,1,12,0,0,0),
    @(173, "Granada", 12,0,0,12,2,0,0),
    @(174, "Fiyi", 12,5,0,12,0,0,0),
    @(175, "Groenlandia", 11,1,3,8,0,0,0),
    @(176, "Curazao", 11,0,3,7,0,0,1),
    @(177, "Seychelles", 10,0,0,10,0,0,0),
    @(178, "Laos", 10,0,0,10,0,0,0),
    @(179, "Mozambique", 10,0,1,9,0,0,0),
    @(180, "Surinam", 10,0,0,9,0,0,1),
    @(181, "Liberia", 10,3,0,9,0,1,1),
    @(182, "Sudan", 10,0,2,6,0,0,2),
    @(183, "Republica del Chad", 9,1,0,9,0,0,0),
    @(184, "Suazilandia", 9,0,0,9,0,0,0),
    @(185, "San Cristobal y Nieves", 9,0,0,9,0,0,0),
    @(186, "Zimbabue", 9,0,0,8,0,0,1),
    @(187, "Nepal", 9,3,1,8,0,0,0),
    @(188, "Montserrat", 9,0,0,7,0,0,2),
    @(189, "Republica de Africa Central", 8,0,0,8,0,0,0),
    @(190, "Angola", 8,0,1,5,0,0,2),
    @(191, "Santa Sede", 7,0,0,7,0,0,0),
    @(192, "San Vicente y las Granadinas", 7,0,1,6,0,0,0),
    @(193, "Somalia", 7,0,1,6,0,0,0),
    @(194, "Cabo Verde", 7,1,0,6,0,0,1),
    @(195, "San Bartolome", 6,0,1,5,0,0,0),
    @(196, "Mauritania", 6,0,2,3,0,0,1),
    @(197, "Islas Turcas y Caicos", 5,0,0,5,0,0,0),
    @(198, "Nicaragua", 5,0,0,4,0,0,1),
    @(199, "Butan", 5,0,2,3,0,0,0),
    @(200, "Sierra Leona", 4,2,0,4,0,0,0),
    @(201, "Belice", 4,0,0,4,0,0,0),
    @(202, "Malaui", 4,1,0,4,0,0,0),
    @(203, "Botsuana", 4,0,0,3,0,0,1),
    @(204, "Gambia", 4,0,2,1,0,0,1),
    @(205, "Burundi", 3,0,0,3,0,0,0),
    @(206, "Anguila", 3,0,0,3,0,0,0),
    @(207, "Islas Virgenes Britanicas", 3,0,0,3,0,0,0),
    @(208, "Bonaire, San Eustaquio y Saba", 2,0,0,2,0,0,0),
    @(209, "Papua Nueva Guinea", 1,0,0,1,0,0,0),
    @(210, "Timor Oriental", 1,0,0,1,0,0,0),
    @(211, "Islas Malvinas", 1,0,0,1,0,0,0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    for ($c = 2; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}

# Update the "last refreshed" footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 20:22"
